$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.924.04"
$ws.Range("E2").Value = "  +6.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.058.36"
$ws.Range("E3").Value = "  +5.93%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.50"
$ws.Range("E5").Value = "  +5.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.32"
$ws.Range("E6").Value = "  +10.45%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.052.80"
$ws.Range("E8").Value = "  +5.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  +7.43%  "

$ws.Range("E10").Value = "  +10.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -3.53%  "

$ws.Range("E12").Value = "  +12.10%  "

$ws.Range("E13").Value = "  +9.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.24"
$ws.Range("E14").Value = "  +9.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.553.52"
$ws.Range("E15").Value = "  +5.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.921.34"
$ws.Range("E16").Value = "  +6.16%  "

$ws.Range("E17").Value = "  +4.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.057.12"
$ws.Range("E18").Value = "  +5.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  +6.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.60"
$ws.Range("E20").Value = "  +6.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  +8.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  +8.53%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.39"
$ws.Range("E23").Value = "  +19.61%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.58"
$ws.Range("E24").Value = "  +11.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.38"
$ws.Range("E25").Value = "  +6.19%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.80"
$ws.Range("E27").Value = "  +7.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("E28").Value = "  +10.14%  "

$ws.Range("E29").Value = "  +7.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.24"
$ws.Range("E31").Value = "  +7.94%  "

$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  +9.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("E34").Value = "  +6.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  +11.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.95"
$ws.Range("E36").Value = "  +3.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0408"
$ws.Range("E37").Value = "  +9.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "445.74"
$ws.Range("E38").Value = "  +3.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0812"
$ws.Range("E39").Value = "  +5.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +25.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.965.27"
$ws.Range("E41").Value = "  +3.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.25"
$ws.Range("E42").Value = "  +7.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.70"
$ws.Range("E44").Value = "  +9.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("E45").Value = "  +11.26%  "

$ws.Range("E46").Value = "  +14.64%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("E48").Value = "  +8.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0516"
$ws.Range("E49").Value = "  +10.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.27"
$ws.Range("E50").Value = "  +5.50%  "

$ws.Range("E51").Value = "  +9.59%  "

